# Insert a new weekly price-report row for "Poroto verde" (Terminal La
# Palmera de La Serena) at row 81, shifting all existing rows 81..179 down
# to 82..180 (dimension grows from A1:R179 to A1:R180).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 81..179 down to 82..180, carrying formats along (this also
# mirrors what Excel's own "Insert Sheet Rows" does on a selected row).
$ws.Rows.Item(81).Insert()

# Populate the freshly inserted row 81 with the new record.
$ws.Range("A81").Value = 8
$ws.Range("B81").Value = "Terminal La Palmera de La Serena"
$ws.Range("C81").Value = "Coquimbo"
$ws.Range("D81").Value = 44589
$ws.Range("E81").Value = 4
$ws.Range("F81").Value = 100112031
$ws.Range("G81").Value = "Poroto verde"
$ws.Range("H81").Value = "Magnum"
$ws.Range("I81").Value = "Primera"
$ws.Range("J81").Value = 400
$ws.Range("K81").Value = 27000
$ws.Range("L81").Value = 28000
$ws.Range("M81").Value = 27500
$ws.Range("N81").Value = "`$/malla 25 kilos"
$ws.Range("O81").Value = "Provincia de Limarí"
$ws.Range("P81").Value = 1100
$ws.Range("Q81").Value = 25
$ws.Range("R81").Value = "Hortaliza"
